$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.575.44"
$ws.Range("E2").Value = "  -0.56%  "
$ws.Range("D3").Value = "'3.120.46"
$ws.Range("E3").Value = "  -1.05%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'598.74"
$ws.Range("E5").Value = "  -2.10%  "
$ws.Range("D6").Value = "'143.27"
$ws.Range("E6").Value = "  -2.75%  "
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("D8").Value = "'3.117.61"
$ws.Range("E8").Value = "  -1.39%  "
$ws.Range("D9").Value = "'0.522"
$ws.Range("E9").Value = "  -0.50%  "
$ws.Range("E10").Value = "  -1.79%  "
$ws.Range("E11").Value = "  -1.06%  "
$ws.Range("E12").Value = "  -1.13%  "
$ws.Range("E13").Value = "  -1.24%  "
$ws.Range("D14").Value = "'35.27"
$ws.Range("E14").Value = "  -0.91%  "
$ws.Range("D15").Value = "'3.632.35"
$ws.Range("E15").Value = "  -1.01%  "
$ws.Range("E16").Value = "  +2.30%  "
$ws.Range("D17").Value = "'63.647.08"
$ws.Range("E17").Value = "  -0.47%  "
$ws.Range("D18").Value = "'3.109.00"
$ws.Range("E18").Value = "  -1.30%  "
$ws.Range("E19").Value = "  -1.29%  "
$ws.Range("D20").Value = "'482.29"
$ws.Range("E20").Value = "  +0.79%  "
$ws.Range("D21").Value = "'14.69"
$ws.Range("E21").Value = "  +0.09%  "
$ws.Range("D22").Value = "'0.708"
$ws.Range("E22").Value = "  -1.08%  "
$ws.Range("D23").Value = "'7.61"
$ws.Range("E23").Value = "  -5.32%  "
$ws.Range("D24").Value = "'87.07"
$ws.Range("E24").Value = "  +4.09%  "
$ws.Range("D25").Value = "'13.34"
$ws.Range("E25").Value = "  -3.26%  "
$ws.Range("E27").Value = "  -2.72%  "
$ws.Range("D28").Value = "'8.24"
$ws.Range("E28").Value = "  -3.67%  "
$ws.Range("D29").Value = "'7.01"
$ws.Range("E29").Value = "  -1.55%  "
$ws.Range("E30").Value = "  -3.34%  "
$ws.Range("D31").Value = "'27.25"
$ws.Range("E31").Value = "  +3.47%  "
$ws.Range("E32").Value = "  -7.34%  "
$ws.Range("E33").Value = "  -0.08%  "
$ws.Range("E34").Value = "  -2.25%  "
$ws.Range("E35").Value = "  -2.65%  "
$ws.Range("E36").Value = "  -0.20%  "
$ws.Range("B37").Value = "PEPE"
$ws.Range("C37").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D37").Value = "'0.0₃0750"
$ws.Range("E37").Value = "  -5.03%  "
$ws.Range("B38").Value = "OKB"
$ws.Range("C38").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D38").Value = "'52.57"
$ws.Range("E38").Value = "  -0.94%  "
$ws.Range("E39").Value = "  -3.73%  "
$ws.Range("D40").Value = "'438.45"
$ws.Range("E40").Value = "  -4.80%  "
$ws.Range("D41").Value = "'0.0394"
$ws.Range("E41").Value = "  -1.14%  "
$ws.Range("E42").Value = "  -0.23%  "
$ws.Range("E43").Value = "  -0.77%  "
$ws.Range("D44").Value = "'2.864.61"
$ws.Range("E44").Value = "  -0.16%  "
$ws.Range("E45").Value = "  -3.38%  "
$ws.Range("E46").Value = "  -4.51%  "
$ws.Range("E47").Value = "  +0.90%  "
$ws.Range("E48").Value = "  +0.04%  "
$ws.Range("D49").Value = "'25.86"
$ws.Range("E49").Value = "  -2.47%  "
$ws.Range("D51").Value = "'121.28"
$ws.Range("E51").Value = "  +2.06%  "
